$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old roster data (rows 6-8 and the header row's H:L empties) ---
$ws.Range("H5:L5").Clear() | Out-Null
$ws.Range("E6:O8").Clear() | Out-Null

# --- Drop the old header text that's being replaced by blanks ---
$ws.Range("E5").ClearContents() | Out-Null
$ws.Range("G5").ClearContents() | Out-Null

# --- Remove the lone hyperlink-styled cell; its style is deleted below ---
$ws.Range("J13").Clear() | Out-Null

# --- New header / data cells (order matters: fixes shared-string indices) ---
$ws.Range("B5").Value = "asdf"
$ws.Range("C5").Value = "adf"
$ws.Range("E4").Value = "test3"
$ws.Range("C4").Value = "test"
$ws.Range("B4").Value = "asdfba"
$ws.Range("D5").Value = "asd"
$ws.Range("G4").Value = "atadd"
$ws.Range("F5").Value = "asdf"

# --- Apply the shared "left aligned" style to every cell that carries it ---
$ws.Range("B4:G4").HorizontalAlignment = -4131
$ws.Range("B5:G5").HorizontalAlignment = -4131
$ws.Range("P10").HorizontalAlignment = -4131
$ws.Range("P11").HorizontalAlignment = -4131
$ws.Range("J12:L12").HorizontalAlignment = -4131
$ws.Range("P12").HorizontalAlignment = -4131

# --- Update the selection to match ---
$ws.Range("B4:G5").Select() | Out-Null

# --- The "Hyperlink" built-in cell style is no longer used anywhere ---
$wb.Styles("Hyperlink").Delete() | Out-Null
